$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column indices: D=4, J=10, K=11, L=12, M=13, P=16
$cols = @(4, 10, 11, 12, 13, 16)

# 1) Snapshot the current ("before") values for rows 20..79 in the columns
#    that move, so the later writes (which happen row-by-row top-to-bottom)
#    never read already-overwritten data.
$old = @{}
for ($r = 20; $r -le 79; $r++) {
    foreach ($c in $cols) {
        $old["$r,$c"] = $ws.Cells.Item($r, $c).Value2
    }
}

# 2) Row 20 keeps its own Volumen/Precio values (J,K,L,M,P) but its Fecha
#    (D) becomes a brand-new data point.
$ws.Cells.Item(20, 4).Value = 44414

# 3) Rows 21..79: each row's D/J/K/L/M/P become the values that used to sit
#    one row above (i.e. the table "grows" by inserting the new row 20
#    reading and pushing everything else down one slot, ending with a new
#    row 80 at the tail).
for ($r = 79; $r -ge 21; $r--) {
    foreach ($c in $cols) {
        $ws.Cells.Item($r, $c).Value = $old["$($r-1),$c"]
    }
}

# 4) Append new row 80, cloning the static columns from (old) row 79 and
#    giving it the shifted-down D/J/K/L/M/P values that used to belong to
#    row 79.
$ws.Cells.Item(80, 1).Value = $ws.Cells.Item(79, 1).Value2
$ws.Cells.Item(80, 2).Value = $ws.Cells.Item(79, 2).Value2
$ws.Cells.Item(80, 3).Value = $ws.Cells.Item(79, 3).Value2
$ws.Cells.Item(80, 4).Value = $old["79,4"]
$ws.Cells.Item(80, 5).Value = $ws.Cells.Item(79, 5).Value2
$ws.Cells.Item(80, 6).Value = $ws.Cells.Item(79, 6).Value2
$ws.Cells.Item(80, 7).Value = $ws.Cells.Item(79, 7).Value2
$ws.Cells.Item(80, 8).Value = $ws.Cells.Item(79, 8).Value2
$ws.Cells.Item(80, 9).Value = $ws.Cells.Item(79, 9).Value2
$ws.Cells.Item(80, 10).Value = $old["79,10"]
$ws.Cells.Item(80, 11).Value = $old["79,11"]
$ws.Cells.Item(80, 12).Value = $old["79,12"]
$ws.Cells.Item(80, 13).Value = $old["79,13"]
$ws.Cells.Item(80, 14).Value = $ws.Cells.Item(79, 14).Value2
$ws.Cells.Item(80, 15).Value = $ws.Cells.Item(79, 15).Value2
$ws.Cells.Item(80, 16).Value = $old["79,16"]
$ws.Cells.Item(80, 17).Value = $ws.Cells.Item(79, 17).Value2
$ws.Cells.Item(80, 18).Value = $ws.Cells.Item(79, 18).Value2

# Match the date number-format style used by the rest of column D.
$ws.Cells.Item(80, 4).NumberFormat = $ws.Cells.Item(79, 4).NumberFormat
